$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection from E6 to E8
$null = $ws.Range("E8").Select()

# Refresh rows 2-18: only the "last refreshed" date (column D) changes
$refreshDate = 45970.394756944443
foreach ($r in 2..18) {
    $ws.Cells.Item($r, 4).Value = $refreshDate
}

# Rows 19-53: terminal-status table refreshed with a new snapshot
# (site name, terminal name, first-unused-since, last-refreshed)
$rowData = @(
    ,@(19, "长沙特来电飞狐四方坪东区充电站", "004A号直流", 45964.528668981482)
    ,@(20, "长沙特来电飞狐四方坪西区充电站", "603号直流", 45966.254062499997)
    ,@(21, "长沙特来电飞狐四方坪南区充电站", "406号直流", 45966.690613425926)
    ,@(22, "长沙特来电飞狐四方坪西区充电站", "505号直流", 45967.507719907408)
    ,@(23, "长沙特来电飞狐四方坪西区充电站", "602号直流", 45967.592800925922)
    ,@(24, "长沙市开福区高岭香江国际城充电站建设项目", "102号直流", 45968.454074074078)
    ,@(25, "长沙特来电飞狐四方坪南区充电站", "502号直流", 45968.550023148149)
    ,@(26, "长沙特来电飞狐四方坪南区充电站", "306号直流", 45969.068564814814)
    ,@(27, "长沙特来电飞狐四方坪东区充电站", "001A号直流", 45969.077048611114)
    ,@(28, "长沙特来电飞狐四方坪西区充电站", "402号直流", 45969.126620370371)
    ,@(29, "长沙市开福区高岭香江国际城充电站建设项目", "109号直流", 45969.429189814815)
    ,@(30, "长沙特来电飞狐四方坪西区充电站", "702号直流", 45969.517708333333)
    ,@(31, "长沙特来电飞狐四方坪西区充电站", "903号直流", 45969.52447916667)
    ,@(32, "长沙特来电飞狐四方坪西区充电站", "B01号直流", 45969.537766203706)
    ,@(33, "长沙特来电飞狐四方坪西区充电站", "404号直流", 45969.542743055557)
    ,@(34, "长沙特来电飞狐四方坪东区充电站", "103号直流", 45969.550613425927)
    ,@(35, "长沙特来电飞狐四方坪东区充电站", "002A号直流", 45969.558425925927)
    ,@(36, "长沙市开福区高岭香江国际城充电站建设项目", "111号直流", 45969.558969907404)
    ,@(37, "长沙特来电飞狐四方坪东区充电站", "001B号直流", 45969.561365740738)
    ,@(38, "长沙特来电飞狐四方坪东区充电站", "005A号直流", 45969.564479166664)
    ,@(39, "长沙市开福区高岭香江国际城充电站建设项目", "203号直流", 45969.566307870373)
    ,@(40, "长沙特来电飞狐四方坪西区充电站", "901号直流", 45969.567974537036)
    ,@(41, "长沙特来电飞狐四方坪南区充电站", "401号直流", 45969.58425925926)
    ,@(42, "长沙市开福区高岭香江国际城充电站建设项目", "210号直流", 45969.585902777777)
    ,@(43, "长沙市开福区高岭香江国际城充电站建设项目", "205号直流", 45969.596851851849)
    ,@(44, "长沙特来电飞狐四方坪西区充电站", "904号直流", 45969.597627314812)
    ,@(45, "长沙特来电飞狐四方坪南区充电站", "304号直流", 45969.608194444445)
    ,@(46, "长沙市开福区高岭香江国际城充电站建设项目", "206号直流", 45969.622754629629)
    ,@(47, "长沙市开福区高岭香江国际城充电站建设项目", "208号直流", 45969.624872685185)
    ,@(48, "长沙特来电飞狐四方坪西区充电站", "A03号直流", 45969.645127314812)
    ,@(49, "长沙特来电飞狐四方坪东区充电站", "006B号直流", 45969.665127314816)
    ,@(50, "长沙特来电飞狐四方坪西区充电站", "501号直流", 45969.672974537039)
    ,@(51, "长沙市开福区高岭香江国际城充电站建设项目", "207号直流", 45969.678668981483)
    ,@(52, "长沙特来电飞狐四方坪东区充电站", "402号直流", 45969.741400462961)
    ,@(53, "长沙市开福区高岭香江国际城充电站建设项目", "204号直流", 45969.771331018521)
)

foreach ($row in $rowData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $refreshDate
}

# Row 54 no longer has data - clear A:D but keep formatting
$null = $ws.Range("A54:D54").ClearContents()
